# Update workbook for data as of 2021-11-07 (October data through 10-30)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-30"

# Update the month label for October in column A, row 11
$ws.Range("A11").Value = "October (through 10-30)"

# Update October row (row 11) values
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = 56
$ws.Range("D11").Value = 79
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 59
$ws.Range("G11").Value = 148
$ws.Range("H11").Value = 189

# Update Total row (row 12) values
$ws.Range("B12").Value = 256
$ws.Range("C12").Value = 485
$ws.Range("D12").Value = 706
$ws.Range("E12").Value = 613
$ws.Range("F12").Value = 481
$ws.Range("G12").Value = 1049
$ws.Range("H12").Value = 1437
